$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 694, shifting existing rows 694..735 down to 695..736
$ws.Rows.Item(694).Insert()

# Populate the newly inserted row 694 with the new data point.
# The date column holds plain text (e.g. "2026/12/29"), not a real date
# value, so force text entry with a leading apostrophe and then strip the
# formatting it leaves behind (quote-prefix style) so the cell ends up
# with no explicit style, matching the rest of the data rows.
$ws.Range("A694").Value = "'2026/01/21"
$ws.Range("A694").ClearFormats()
$ws.Range("B694").Value = "水"
$ws.Range("C694").Value = 6
$ws.Range("D694").Value = 196
